# "Sell_analyse" report: rebuild sheet1 as a small "Эффективность запасов"
# (Inventory efficiency) table with a title row, a header row, a template
# row (Smarty-style placeholders) and a blank bottom row, boxed in borders.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---- values --------------------------------------------------------------
# Row 1: merged title
$ws.Range("B1").Value = "Эффективность запасов"
# Row 2: column headers
$ws.Range("B2").Value = "Категория"
$ws.Range("C2").Value = "Групировка"
$ws.Range("D2").Value = "Среднее"
# Row 3: Smarty template placeholders (the actual report row)
$ws.Range("B3").Value = "{`$v->rows[]->cathegory}"
$ws.Range("C3").Value = "{`$v->rows[]->group_by}"
$ws.Range("D3").Value = "{`$v->rows[]->avg}"
# Row 4 (B4:D4) intentionally left blank - bottom border row only.

# ---- title formatting ------------------------------------------------------
$title = $ws.Range("B1:D1")
$title.Font.Bold = $true
$title.Font.Size = 16
$title.HorizontalAlignment = -4108   # xlCenter
[void]$title.Merge()

# ---- row heights -----------------------------------------------------------
$ws.Rows("1").RowHeight = 21
$ws.Rows("4").RowHeight = 15.75

# ---- borders: thin grid everywhere, medium (thick) box around the outside -
$grid = $ws.Range("B1:D4")
$grid.Borders.LineStyle = 1
$grid.Borders.Weight = 2             # xlThin
[void]$grid.BorderAround(1, -4138)   # xlContinuous, xlMedium

# ---- column widths ----------------------------------------------------------
$ws.Columns("B").ColumnWidth = 21.8
$ws.Columns("C").ColumnWidth = 21.5
$ws.Columns("D").ColumnWidth = 15.8

# ---- sheet view: no gridlines, selection parked at G7 -----------------------
$excel.ActiveWindow.DisplayGridlines = $false
[void]$ws.Range("G7").Select()

# ---- page setup --------------------------------------------------------------
$ws.PageSetup.PaperSize = 9          # xlPaperA4
$ws.PageSetup.Orientation = 1        # xlPortrait
